# CREATE scripts and working hours adjusted
#
# Fill in the previously-empty "Documentation" task row (row 6) on the
# Janeczek worksheet with the date it happened, the task description,
# and the estimated / actual working hours spent on it. Dependent SUM
# formulas (row 12 on this sheet, and the SUM worksheet) recalculate
# automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Janeczek")

$ws.Range("B6").Value = 41956
$ws.Range("C6").Value = "Documentation"
$ws.Range("D6").Value = "Technology Description + Further Approach"
$ws.Range("E6").Value = 0.041666666666666664
$ws.Range("F6").Value = 0.020833333333333332

# Janeczek is now the focused / active sheet, with G6 selected.
$ws.Activate()
$ws.Range("G6").Select()
